# Apply updated market-price values scraped by the scheduled runner.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 291.58334
$ws.Range("I5").Value = 237.375
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 237.375
$ws.Range("L5").Value = 400
$ws.Range("M5").Value = -122.375
$ws.Range("N5").Value = -630
$ws.Range("H106").Value = 11787.538
$ws.Range("I106").Value = 12269.333
$ws.Range("K106").Value = 12269.333
$ws.Range("M106").Value = -11638.333
$ws.Range("H111").Value = 2729.5334
$ws.Range("I111").Value = 2335.6667
$ws.Range("K111").Value = 7007.000100000001
$ws.Range("M111").Value = -3940.000100000001
$ws.Range("H113").Value = 2881
$ws.Range("J113").Value = 3179
$ws.Range("L113").Value = 3179
$ws.Range("N113").Value = -9687
$ws.Range("H132").Value = 6066538
$ws.Range("I132").Value = 8551115
$ws.Range("J132").Value = 10381.9375
$ws.Range("K132").Value = 25653345
$ws.Range("L132").Value = 31145.8125
$ws.Range("M132").Value = -25650815
$ws.Range("N132").Value = -36205.8125
$ws.Range("H135").Value = 27027584
$ws.Range("I135").Value = 230.30302
$ws.Range("K135").Value = 2072.72718
$ws.Range("M135").Value = 462.2728200000001
$ws.Range("H137").Value = 1097.6301
$ws.Range("I137").Value = 894.20514
$ws.Range("J137").Value = 1330.9706
$ws.Range("K137").Value = 2682.61542
$ws.Range("L137").Value = 3992.9118
$ws.Range("M137").Value = -132.6154200000001
$ws.Range("N137").Value = -9092.9118
$ws.Range("H138").Value = 725909
$ws.Range("I138").Value = 927.6667
$ws.Range("J138").Value = 1611997.2
$ws.Range("K138").Value = 2783.0001
$ws.Range("L138").Value = 4835991.6
$ws.Range("M138").Value = 2356.9999
$ws.Range("N138").Value = -4846271.6
$ws.Range("H141").Value = 557.05884
$ws.Range("I141").Value = 557.05884
$ws.Range("K141").Value = 1671.17652
$ws.Range("M141").Value = 3508.82348

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 216.83333
$ws.Range("I5").Value = 290.25
$ws.Range("J5").Value = 70
$ws.Range("K5").Value = 290.25
$ws.Range("L5").Value = 70
$ws.Range("M5").Value = -178.25
$ws.Range("N5").Value = -294
$ws.Range("H32").Value = 4516.8594
$ws.Range("I32").Value = 4126.9106
$ws.Range("J32").Value = 7246.5
$ws.Range("K32").Value = 4126.9106
$ws.Range("L32").Value = 7246.5
$ws.Range("M32").Value = -3839.9106
$ws.Range("N32").Value = -7820.5
$ws.Range("H63").Value = 55557696
$ws.Range("I63").Value = 2020.8182
$ws.Range("K63").Value = 2020.8182
$ws.Range("M63").Value = -1334.8182
$ws.Range("H66").Value = 55557696
$ws.Range("I66").Value = 2020.8182
$ws.Range("K66").Value = 10104.091
$ws.Range("M66").Value = -6672.091
$ws.Range("H74").Value = 1177.1351
$ws.Range("I74").Value = 811.7241
$ws.Range("J74").Value = 2501.75
$ws.Range("K74").Value = 811.7241
$ws.Range("L74").Value = 2501.75
$ws.Range("M74").Value = 62.27589999999998
$ws.Range("N74").Value = -4249.75
$ws.Range("H77").Value = 1177.1351
$ws.Range("I77").Value = 811.7241
$ws.Range("J77").Value = 2501.75
$ws.Range("K77").Value = 4058.6205
$ws.Range("L77").Value = 12508.75
$ws.Range("M77").Value = 309.3795
$ws.Range("N77").Value = -21244.75
$ws.Range("H132").Value = 3178.2173
$ws.Range("I132").Value = 3870.3572
$ws.Range("J132").Value = 2101.5557
$ws.Range("K132").Value = 11611.0716
$ws.Range("L132").Value = 6304.6671
$ws.Range("M132").Value = -9081.071599999999
$ws.Range("N132").Value = -11364.6671

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 216.83333
$ws.Range("I4").Value = 290.25
$ws.Range("J4").Value = 70
$ws.Range("K4").Value = 290.25
$ws.Range("L4").Value = 70
$ws.Range("M4").Value = -175.25
$ws.Range("N4").Value = -300
$ws.Range("H134").Value = 10251.8
$ws.Range("I134").Value = 1972.125
$ws.Range("K134").Value = 5916.375
$ws.Range("M134").Value = -3381.375
$ws.Range("H137").Value = 28623.8
$ws.Range("J137").Value = 28623.8
$ws.Range("L137").Value = 28623.8
$ws.Range("N137").Value = -38823.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 367.22223
$ws.Range("I7").Value = 244
$ws.Range("J7").Value = 521.25
$ws.Range("K7").Value = 244
$ws.Range("L7").Value = 521.25
$ws.Range("M7").Value = -131
$ws.Range("N7").Value = -747.25
$ws.Range("H31").Value = 2051.5557
$ws.Range("I31").Value = 2253.3
$ws.Range("K31").Value = 2253.3
$ws.Range("M31").Value = -1958.3
$ws.Range("H34").Value = 2051.5557
$ws.Range("I34").Value = 2253.3
$ws.Range("K34").Value = 2253.3
$ws.Range("M34").Value = -2051.3
$ws.Range("H58").Value = 850.8125
$ws.Range("I58").Value = 763.4167
$ws.Range("J58").Value = 1113
$ws.Range("K58").Value = 763.4167
$ws.Range("L58").Value = 1113
$ws.Range("M58").Value = -560.4167
$ws.Range("N58").Value = -1519
$ws.Range("H132").Value = 3070.7
$ws.Range("I132").Value = 2705.2632
$ws.Range("J132").Value = 10014
$ws.Range("K132").Value = 8115.7896
$ws.Range("L132").Value = 30042
$ws.Range("M132").Value = -5585.7896
$ws.Range("N132").Value = -35102
$ws.Range("H134").Value = 17858266
$ws.Range("I134").Value = 1125.8823
$ws.Range("J134").Value = 45455664
$ws.Range("K134").Value = 3377.6469
$ws.Range("L134").Value = 136366992
$ws.Range("M134").Value = -842.6468999999997
$ws.Range("N134").Value = -136372062
$ws.Range("H136").Value = 850.8125
$ws.Range("I136").Value = 763.4167
$ws.Range("J136").Value = 1113
$ws.Range("K136").Value = 2290.2501
$ws.Range("L136").Value = 3339
$ws.Range("M136").Value = 259.7498999999998
$ws.Range("N136").Value = -8439

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 90.935486
$ws.Range("J12").Value = 74.478264
$ws.Range("L12").Value = 223.434792
$ws.Range("N12").Value = -569.434792
$ws.Range("H119").Value = 6675.8667
$ws.Range("I119").Value = 809.3333
$ws.Range("J119").Value = 8142.5
$ws.Range("K119").Value = 2427.9999
$ws.Range("L119").Value = 24427.5
$ws.Range("M119").Value = 2410.0001
$ws.Range("N119").Value = -34103.5
$ws.Range("H127").Value = 2497
$ws.Range("J127").Value = 2497
$ws.Range("L127").Value = 7491
$ws.Range("N127").Value = -17411
$ws.Range("H131").Value = 13890049
$ws.Range("J131").Value = 1217.3135
$ws.Range("L131").Value = 3651.9405
$ws.Range("N131").Value = -13731.9405
$ws.Range("H139").Value = 1518.4651
$ws.Range("I139").Value = 1398.0416
$ws.Range("J139").Value = 1670.579
$ws.Range("K139").Value = 4194.1248
$ws.Range("L139").Value = 5011.737
$ws.Range("M139").Value = 945.8752000000004
$ws.Range("N139").Value = -15291.737

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3999.8
$ws.Range("I122").Value = 3333.3333
$ws.Range("K122").Value = 9999.999899999999
$ws.Range("M122").Value = -7549.999899999999
$ws.Range("H132").Value = 1655.9642
$ws.Range("I132").Value = 1473.6666
$ws.Range("J132").Value = 2749.75
$ws.Range("K132").Value = 4420.9998
$ws.Range("L132").Value = 8249.25
$ws.Range("M132").Value = -1890.9998
$ws.Range("N132").Value = -13309.25
$ws.Range("H135").Value = 44999.5
$ws.Range("J135").Value = 39999
$ws.Range("L135").Value = 39999
$ws.Range("N135").Value = -50139
$ws.Range("H136").Value = 26489.9
$ws.Range("J136").Value = 26489.9
$ws.Range("L136").Value = 79469.70000000001
$ws.Range("N136").Value = -84569.70000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1044.1052
$ws.Range("I16").Value = 1149.7333
$ws.Range("J16").Value = 648
$ws.Range("K16").Value = 1149.7333
$ws.Range("L16").Value = 648
$ws.Range("M16").Value = -979.7333000000001
$ws.Range("N16").Value = -988
$ws.Range("H100").Value = 1742.3
$ws.Range("I100").Value = 1654.8
$ws.Range("J100").Value = 1829.8
$ws.Range("K100").Value = 1654.8
$ws.Range("L100").Value = 1829.8
$ws.Range("M100").Value = -1113.8
$ws.Range("N100").Value = -2911.8
$ws.Range("H132").Value = 15574.3
$ws.Range("I132").Value = 970.8444
$ws.Range("J132").Value = 41860.52
$ws.Range("K132").Value = 2912.5332
$ws.Range("L132").Value = 125581.56
$ws.Range("M132").Value = -382.5331999999999
$ws.Range("N132").Value = -130641.56
$ws.Range("H136").Value = 1126.25
$ws.Range("I136").Value = 1051.3334
$ws.Range("K136").Value = 3154.0002
$ws.Range("M136").Value = -604.0001999999999
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 504.2143
$ws.Range("I107").Value = 470
$ws.Range("K107").Value = 1410
$ws.Range("M107").Value = 510
$ws.Range("H108").Value = 23333.334
$ws.Range("J108").Value = 23333.334
$ws.Range("L108").Value = 23333.334
$ws.Range("N108").Value = -31013.334
$ws.Range("H109").Value = 36780.668
$ws.Range("H132").Value = 4348.95
$ws.Range("I132").Value = 5315.2666
$ws.Range("J132").Value = 1450
$ws.Range("K132").Value = 15945.7998
$ws.Range("L132").Value = 4350
$ws.Range("M132").Value = -13415.7998
$ws.Range("N132").Value = -9410
$ws.Range("H136").Value = 618.3684
$ws.Range("I136").Value = 438.86206
$ws.Range("J136").Value = 1196.7778
$ws.Range("K136").Value = 1316.58618
$ws.Range("L136").Value = 3590.3334
$ws.Range("M136").Value = 1233.41382
$ws.Range("N136").Value = -8690.3334
